# "Generate Report for Handback"
#
# The handback transform for the "30a979c7-27d7-4f3b-8f20-2e774d6187e8" file
# failed because the handback file name didn't match the handoff file name.
# Update the localization-status report to reflect this:
#   - Overview sheet: status for that file becomes "Handback transform failed"
#   - zh-cn / de-de sheets: Status column becomes "Handback transform failed"
#     and the (previously empty) Error Detail column is filled in with the
#     mismatch explanation for each locale.
#   - The Error Detail column is widened so the message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the 30a979c7... file, zh-cn/de-de status columns (E/F)
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# zh-cn / de-de detail sheets: row 3 is the same file, Status is column C
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Fill in the Error Detail column (P) with the mismatch message for each locale
$zhcn.Range("P3").Value = "Handback file name: glxfnqpx.q5p is different with handoff file name: 30a979c7-27d7-4f3b-8f20-2e774d6187e8.59cc4c485659c0ed65e1c665cef021bf89b153bd.zh-cn."
$dede.Range("P3").Value = "Handback file name: glxfnqpx.q5p is different with handoff file name: 30a979c7-27d7-4f3b-8f20-2e774d6187e8.59cc4c485659c0ed65e1c665cef021bf89b153bd.de-de."

# Widen the Error Detail column (column 16 / P) on both sheets so the message
# is legible. 39.1666... (character units) is stored by Excel as width="40".
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
